$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input Equations")

# ------------------------------------------------------------------
# 1) Add the "175 gal drop tank" rows, right after the existing
#    "120 gal drop tank" rows, before the "EO/ISR" section.
# ------------------------------------------------------------------
$ws.Rows("156:157").Insert()

$ws.Range("A156").Value = "175 gal drop tank d"
$ws.Range("B156").Value = "d_175dropTank"
$ws.Range("C156").Value = 27
$ws.Range("D156").Value = "in"

$ws.Range("A157").Value = "175 gal drop tank"
$ws.Range("B157").Value = "area175DropTank"
$ws.Range("C157").Formula = "=((C156^2 *PI())/4)/144"
$ws.Range("D157").Value = "ft^2"

# ------------------------------------------------------------------
# 2) Remove two redundant blank spacer rows above the
#    "Drop Tank Corrections" header.
# ------------------------------------------------------------------
$ws.Rows("160:161").Delete()
$ws.Rows("162:162").Delete()

# ------------------------------------------------------------------
# 3) Add one extra blank spacer row right below the row that carries
#    the F-column note.
# ------------------------------------------------------------------
$ws.Rows("168:168").Insert()
$ws.Range("F168").Delete()

# ------------------------------------------------------------------
# 4) Remove the now-redundant blank spacer row above the
#    "Fuel Jettison" header.
# ------------------------------------------------------------------
$ws.Rows("172:172").Delete()
